$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" ---
# This string is shared across all three sheets (Overview, zh-cn, de-de).
# Re-assigning every cell that currently holds it collapses back down to a
# single shared string entry once the old text is no longer referenced.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column widths: narrow the "handoff/handback" datetime columns ---
# Overview sheet: columns E and F
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
